# "removed duplicate 128-11 run"
#
# The Data sheet's last row (row 28, date 2016-05-15 / serial 42501) had one
# extra PTC run counted twice. Removing the duplicate drops:
#   - Total PTC Runs   (B28): 145 -> 144
#   - Total Completed  (C28): 141 -> 140
#   - Cut Out Runs     (F28): 141 -> 140
# which in turn changes the (statically cached, not formula-driven) derived
# figures for that row:
#   - Total Completed %            (G28): 0.97241379310344822 -> 0.97222222222222221
#   - Completed Trip Length Average(H28): 43.178160919117389  -> 43.391666666163864
#
# Single/Multiple Init Runs (D28/E28) and the trip-length min/max (I28/J28)
# are untouched by the fix.
#
# The author also left the selection on the corrected row (H28:J28) when the
# workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B28").Value = 144
$ws.Range("C28").Value = 140
$ws.Range("F28").Value = 140
$ws.Range("G28").Value = 0.97222222222222221
$ws.Range("H28").Value = 43.391666666163864

$ws.Range("H28:J28").Select() | Out-Null
